$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value2 = "62.760.24"
$ws.Range("E2").Value2 = "  +6.36%  "

# Row 3 (Ethereum)
$ws.Range("D3").Value2 = "3.108.95"
$ws.Range("E3").Value2 = "  +3.74%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value2 = "  +0.00%  "

# Row 5 (BNB)
$ws.Range("D5").Value2 = "584.90"
$ws.Range("E5").Value2 = "  +4.39%  "

# Row 6 (Solana)
$ws.Range("D6").Value2 = "143.44"
$ws.Range("E6").Value2 = "  +4.69%  "

# Row 7 (USDC)
$ws.Range("E7").Value2 = "  -0.07%  "

# Row 8 (LidoStakedEther)
$ws.Range("D8").Value2 = "3.101.39"
$ws.Range("E8").Value2 = "  +3.77%  "

# Row 9 (XRP)
$ws.Range("D9").Value2 = "0.532"
$ws.Range("E9").Value2 = "  +2.02%  "

# Row 10 (Dogecoin)
$ws.Range("E10").Value2 = "  +9.00%  "

# Row 11 (Toncoin)
$ws.Range("E11").Value2 = "  +9.84%  "

# Row 12 (Cardano)
$ws.Range("D12").Value2 = "0.469"
$ws.Range("E12").Value2 = "  +2.88%  "

# Row 13 (ShibaInu)
$ws.Range("D13").Value2 = "0.0000243"
$ws.Range("E13").Value2 = "  +5.87%  "

# Row 14 (Avalanche)
$ws.Range("D14").Value2 = "35.57"
$ws.Range("E14").Value2 = "  +5.94%  "

# Row 15 (TRON)
$ws.Range("E15").Value2 = "  +0.90%  "

# Row 16 (WrappedliquidstakedEther2.0)
$ws.Range("D16").Value2 = "3.626.58"
$ws.Range("E16").Value2 = "  +3.75%  "

# Row 17 (Polkadot)
$ws.Range("D17").Value2 = "7.29"
$ws.Range("E17").Value2 = "  +0.16%  "

# Row 18 (WrappedEther)
$ws.Range("D18").Value2 = "3.112.76"
$ws.Range("E18").Value2 = "  +3.91%  "

# Row 19 (WrappedBTC)
$ws.Range("D19").Value2 = "62.694.71"
$ws.Range("E19").Value2 = "  +6.20%  "

# Row 20 (BitcoinCash)
$ws.Range("D20").Value2 = "453.93"
$ws.Range("E20").Value2 = "  +6.12%  "

# Row 21 (Chainlink)
$ws.Range("D21").Value2 = "14.06"
$ws.Range("E21").Value2 = "  +2.53%  "

# Row 22 (Polygon)
$ws.Range("D22").Value2 = "0.735"

# Row 23 (Uniswap)
$ws.Range("D23").Value2 = "7.54"
$ws.Range("E23").Value2 = "  +6.12%  "

# Row 24 (InternetComputer(DFINITY))
$ws.Range("D24").Value2 = "13.70"
$ws.Range("E24").Value2 = "  +2.74%  "

# Row 25 (Litecoin)
$ws.Range("D25").Value2 = "82.12"
$ws.Range("E25").Value2 = "  +1.82%  "

# Row 26 (Dai)
$ws.Range("D26").Value2 = "1.00"
$ws.Range("E26").Value2 = "  -0.07%  "

# Row 27 (ImmutableX)
$ws.Range("D27").Value2 = "2.28"
$ws.Range("E27").Value2 = "  +4.12%  "

# Row 28 (PancakeSwap)
$ws.Range("D28").Value2 = "2.69"
$ws.Range("E28").Value2 = "  +6.01%  "

# Row 29 (FirstDigitalUSD)
$ws.Range("E29").Value2 = "  +0.01%  "

# Row 30 (RenderToken)
$ws.Range("D30").Value2 = "8.24"
$ws.Range("E30").Value2 = "  +5.70%  "

# Row 31 (NEARProtocol)
$ws.Range("D31").Value2 = "6.88"
$ws.Range("E31").Value2 = "  +13.94%  "

# Row 32 (Hedera)
$ws.Range("E32").Value2 = "  +12.89%  "

# Row 33 (EthereumClassic)
$ws.Range("D33").Value2 = "27.10"
$ws.Range("E33").Value2 = "  +5.32%  "

# Row 34 (Mantle)
$ws.Range("E34").Value2 = "  +4.68%  "

# Row 35 (PEPE)
$ws.Range("D35").Value2 = "0.0₃0802"
$ws.Range("E35").Value2 = "  +5.94%  "

# Row 36 (Filecoin)
$ws.Range("D36").Value2 = "6.09"
$ws.Range("E36").Value2 = "  +2.51%  "

# Row 37 (Stacks)
$ws.Range("E37").Value2 = "  +6.32%  "

# Row 38 <-> Row 39 swap: dogwifhat/OKB rows exchange places (with slightly different values)
$ws.Range("B38").Value2 = "OKB"
$ws.Range("C38").Value2 = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").Value2 = "50.71"
$ws.Range("E38").Value2 = "  +4.09%  "

$ws.Range("B39").Value2 = "dogwifhat"
$ws.Range("C39").Value2 = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").Value2 = "3.04"
$ws.Range("E39").Value2 = "  +11.32%  "

# Row 40 (Cosmos)
$ws.Range("D40").Value2 = "8.80"
$ws.Range("E40").Value2 = "  +1.58%  "

# Row 41 (Bittensor)
$ws.Range("D41").Value2 = "427.00"
$ws.Range("E41").Value2 = "  +7.04%  "

# Row 42 (Maker)
$ws.Range("D42").Value2 = "2.932.28"
$ws.Range("E42").Value2 = "  +6.44%  "

# Row 43 (VeChain)
$ws.Range("E43").Value2 = "  +6.12%  "

# Row 44 (TheGraph)
$ws.Range("D44").Value2 = "0.282"
$ws.Range("E44").Value2 = "  +12.44%  "

# Row 45 (Kaspa)
$ws.Range("D45").Value2 = "0.112"
$ws.Range("E45").Value2 = "  +3.60%  "

# Row 46 (Fetch.AI)
$ws.Range("D46").Value2 = "2.17"
$ws.Range("E46").Value2 = "  +8.67%  "

# Row 47 (Monero)
$ws.Range("D47").Value2 = "125.24"
$ws.Range("E47").Value2 = "  +1.51%  "

# Row 48 <-> Row 49 swap: Arweave/USDe rows exchange places (with slightly different values)
$ws.Range("B48").Value2 = "USDe"
$ws.Range("C48").Value2 = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").Value2 = "0.999"
$ws.Range("E48").Value2 = "  -0.01%  "

$ws.Range("B49").Value2 = "Arweave"
$ws.Range("C49").Value2 = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D49").Value2 = "35.01"
$ws.Range("E49").Value2 = "  -1.07%  "

# Row 50 (Stellar)
$ws.Range("E50").Value2 = "  +1.40%  "

# Row 51 (InjectiveProtocol)
$ws.Range("D51").Value2 = "24.65"
$ws.Range("E51").Value2 = "  +5.37%  "
